# Apply the "Fixed bugs caused by restructuring for git" edit:
# add hidden helper lists (Color / Line / Marker) in columns I:K,
# wire them up as data-validation source lists for columns E:G,
# and fix the E2 header + C3 default value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing cells -----------------------------------------------
$ws.Range("E2").Value = "color (opt)"
$ws.Range("C3").Value = 0

# Populate the visible E3:G3 selections to match the new helper lists
$ws.Range("E3").Value = "k"
$ws.Range("F3").Value = "-."
$ws.Range("G3").Value = "<"

# --- Helper header row (I2:K2), styled like the other header cells ----
$ws.Range("I2").Value = "Color"
$ws.Range("J2").Value = "Line"
$ws.Range("K2").Value = "Marker"
$ws.Range("B2").Copy()
$ws.Range("I2:K2").PasteSpecial(-4122)

# --- Color helper list (I3:I8) -----------------------------------------
$ws.Range("I3").Value = "k"
$ws.Range("I4").Value = "r"
$ws.Range("I5").Value = "b"
$ws.Range("I6").Value = "g"
$ws.Range("I7").Value = "m"
$ws.Range("I8").Value = "c"

# --- Line-style helper list (J3:J6) ------------------------------------
$ws.Range("J3").Value = "-"
$ws.Range("J4").Value = "'--"
$ws.Range("J5").Value = "'-."
$ws.Range("J6").Value = ":"

# --- Marker helper list (K3:K11) ---------------------------------------
$ws.Range("K3").Value = "o"
$ws.Range("K4").Value = "s"
$ws.Range("K5").Value = "d"
$ws.Range("K6").Value = "<"
$ws.Range("K7").Value = "^"
$ws.Range("K8").Value = ">"
$ws.Range("K9").Value = "x"
$ws.Range("K10").Value = "."
$ws.Range("K11").Value = "*"

# --- small purple legend-style cell next to the new helper table -------
$ws.Range("L4").Font.Size = 10
$ws.Range("L4").Font.Color = 15736992
$ws.Range("L4").Font.Name = "Helvetica"

# --- Hide the helper columns -------------------------------------------
$helperCols = $ws.Range("I1:K1").EntireColumn
$helperCols.ColumnWidth = 0
$helperCols.Hidden = $true

# --- Wire up data validation lists for E, F, G using the helper ranges -
$ws.Range("E3:E12").Validation.Add(3, 1, 1, "=`$I`$3:`$I`$8")
$ws.Range("F3:F12").Validation.Add(3, 1, 1, "=`$J`$3:`$J`$6")
$ws.Range("G3:G12").Validation.Add(3, 1, 1, "=`$K`$3:`$K`$11")
